$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the new registry entry on row 50.
$ws.Range("A50").Value = "ISO_"
$ws.Range("D50").Value = "betsy.fanning@3dpdfconsortium.com"
$ws.Range("E50").Value = "ISO (via the 3D PDF Consortium)"
$ws.Range("F50").Value = 43614

# Turn the email address into a mailto hyperlink, then restore the shared
# "Hyperlink" cell style (Add() applies its own ad-hoc style otherwise).
$ws.Hyperlinks.Add($ws.Range("D50"), "mailto:betsy.fanning@3dpdfconsortium.com")
$ws.Range("D50").Style = "Hyperlink"

# Move the active selection from F50 to B50.
$ws.Range("B50").Select()
